# "aula do dia 25" - add a summary-statistics block (mean/median/mode/
# quartiles/percentiles/min/max of salario) to the BaseDados sheet and
# switch the active tab back to BaseDados.

$wb  = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("BaseDados")
$ws2 = $wb.Worksheets.Item("Código R")

# --- BaseDados: add the summary-statistics block in columns H:I (rows 3-13) ---
# Labels/formulas are entered in the same order the author originally typed
# them (H7 / "q2 / Mediana" last), so newly-interned shared strings land in
# the same order as in the saved workbook.
$ws1.Range("H3").Value = "Média"
$ws1.Range("I3").Formula = '=AVERAGE(F2:F37)'

$ws1.Range("H4").Value = "Mediana"
$ws1.Range("I4").Formula = '=MEDIAN(F1:F37)'

$ws1.Range("H5").Value = "Moda"
$ws1.Range("I5").Formula = '=MODE(F2:F37)'

$ws1.Range("H6").Value = "q1"
$ws1.Range("I6").Formula = '=QUARTILE.INC(F$2:F$37,1)'

$ws1.Range("H8").Value = "q3"
$ws1.Range("I8").Formula = '=QUARTILE.INC(F$2:F$37,3)'

$ws1.Range("H9").Value = "q4"
$ws1.Range("I9").Formula = '=QUARTILE.INC(F$2:F$37,4)'

$ws1.Range("H10").Value = "p1"
$ws1.Range("I10").Formula = '=PERCENTILE(F$2:F$37,0.01)'

$ws1.Range("H11").Value = "p2"
$ws1.Range("I11").Formula = '=PERCENTILE(F$2:F$37,0.99)'

$ws1.Range("H12").Value = "Minimo"
$ws1.Range("I12").Formula = '=MIN(F2:F37)'

$ws1.Range("H13").Value = "Maximo"
$ws1.Range("I13").Formula = '=MAX(F2:F37)'

$ws1.Range("H7").Value = "q2 / Mediana"
$ws1.Range("I7").Formula = '=QUARTILE.INC(F$2:F$37,2)'

# New column H needs to be sized to fit its new label text.
$ws1.Columns.Item(8).ColumnWidth = 12.5703125

# --- Active sheet / selection moves from "Código R" back to "BaseDados" ---
# (sheet2's own selection, A4, was already in place and must stay untouched)
$ws1.Activate()
$ws1.Range("H7").Select()
